# Apply odds updates to "Jogos da Semana" worksheet (match week 2025-04-03)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 4.33
$ws.Range("L2").Value = 1.14
$ws.Range("M2").Value = 5.5
$ws.Range("P2").Value = 1.25
$ws.Range("Q2").Value = 3.75
$ws.Range("R2").Value = 1.47
$ws.Range("S2").Value = 2.5
$ws.Range("Y2").Value = 21
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 9.5
$ws.Range("AB2").Value = 15
$ws.Range("AC2").Value = 41
$ws.Range("AD2").Value = 126
$ws.Range("AE2").Value = 21
$ws.Range("AJ2").Value = 34

# Row 3
$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 1.05
$ws.Range("K3").Value = 11
$ws.Range("L3").Value = 1.29
$ws.Range("M3").Value = 3.5
$ws.Range("N3").Value = 1.95
$ws.Range("O3").Value = 1.85
$ws.Range("R3").Value = 1.83
$ws.Range("S3").Value = 1.83
$ws.Range("T3").Value = 7
$ws.Range("U3").Value = 8.5
$ws.Range("V3").Value = 8.5
$ws.Range("W3").Value = 15
$ws.Range("X3").Value = 15
$ws.Range("AB3").Value = 15
$ws.Range("AC3").Value = 51
$ws.Range("AD3").Value = 251
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 23
$ws.Range("AG3").Value = 15
$ws.Range("AH3").Value = 51
$ws.Range("AI3").Value = 41
$ws.Range("AJ3").Value = 41

# Row 4
$ws.Range("G4").Value = 2.62
$ws.Range("H4").Value = 2.7
$ws.Range("I4").Value = 3
$ws.Range("T4").Value = 7.6
$ws.Range("U4").Value = 13.5
$ws.Range("V4").Value = 9.5
$ws.Range("W4").Value = 32
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 32
$ws.Range("AE4").Value = 8.25
$ws.Range("AF4").Value = 15.5
$ws.Range("AG4").Value = 10.25
$ws.Range("AH4").Value = 40
$ws.Range("AI4").Value = 28
$ws.Range("AJ4").Value = 35

# Row 5
$ws.Range("T5").Value = 16.5
$ws.Range("X5").Value = 50
$ws.Range("Y5").Value = 45

# Row 6
$ws.Range("G6").Value = 1.38
$ws.Range("I6").Value = 7.5
$ws.Range("J6").Value = 1.03
$ws.Range("K6").Value = 17
$ws.Range("P6").Value = 1.25
$ws.Range("Q6").Value = 3.75
$ws.Range("T6").Value = 9.5
$ws.Range("AC6").Value = 41
$ws.Range("AD6").Value = 151
$ws.Range("AE6").Value = 23

# Row 7
$ws.Range("T7").Value = 11

